$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '38.771.38'
$ws.Range("E2").Value = '  -0.25%  '

$ws.Range("D3").Value = '2.100.50'
$ws.Range("E3").Value = '  -0.22%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.33'
$ws.Range("E5").Value = '  -0.60%  '

$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.74'
$ws.Range("E7").Value = '  +2.16%  '

$ws.Range("E9").Value = '  +0.78%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0841'
$ws.Range("E10").Value = '  +0.46%  '

$ws.Range("E11").Value = '  -1.15%  '

$ws.Range("E12").Value = '  +5.21%  '

$ws.Range("D13").Value = '2.410.78'
$ws.Range("E13").Value = '  -0.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.99'
$ws.Range("E14").Value = '  -1.00%  '

$ws.Range("E15").Value = '  +0.37%  '

$ws.Range("E16").Value = '  +0.04%  '

$ws.Range("D17").Value = '2.088.56'
$ws.Range("E17").Value = '  -0.85%  '

$ws.Range("D18").Value = '38.752.23'
$ws.Range("E18").Value = '  -0.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.59'
$ws.Range("E19").Value = '  -0.70%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.02'
$ws.Range("E20").Value = '  -0.51%  '

$ws.Range("E21").Value = '  +0.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.83'
$ws.Range("E22").Value = '  +0.18%  '

$ws.Range("E24").Value = '  -3.92%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.31'
$ws.Range("E25").Value = '  -1.42%  '

$ws.Range("E26").Value = '  +1.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.00'
$ws.Range("E27").Value = '  -0.44%  '

$ws.Range("E28").Value = '  -0.60%  '

$ws.Range("E29").Value = '  +1.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.32'
$ws.Range("E30").Value = '  +0.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.54'
$ws.Range("E31").Value = '  +8.82%  '

$ws.Range("E32").Value = '  +0.18%  '

$ws.Range("B34").Value = 'THORChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.17'
$ws.Range("E34").Value = '  +11.66%  '

$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.79'
$ws.Range("E35").Value = '  +0.71%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0613'
$ws.Range("E36").Value = '  -0.14%  '

$ws.Range("E37").Value = '  -1.77%  '

$ws.Range("E38").Value = '  -0.63%  '

$ws.Range("E39").Value = '  +0.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.97'
$ws.Range("E40").Value = '  -2.33%  '

$ws.Range("E41").Value = '  +3.42%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.74'
$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("D43").Value = '1.525.30'
$ws.Range("E43").Value = '  -1.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.20'
$ws.Range("E44").Value = '  +7.18%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.81'
$ws.Range("E45").Value = '  -0.50%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.75'
$ws.Range("E46").Value = '  +0.64%  '

$ws.Range("E47").Value = '  -2.22%  '

$ws.Range("E48").Value = '  +4.30%  '

$ws.Range("E49").Value = '  +0.93%  '

$ws.Range("E50").Value = '  -1.23%  '

$ws.Range("D51").Value = '2.298.12'
$ws.Range("E51").Value = '  -0.30%  '
